$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New task rows 33-36: fill in the activity text / dates / trailing style cells ---
# (rows 33-36 already existed with only a B-column index number; we add the rest)

# Copy date-cell formatting (style index 6) from an existing date cell (D32/E32) onto the new ones
$ws.Range("D32").Copy()
$ws.Range("D33:D39").PasteSpecial(-4122)
$ws.Range("E32").Copy()
$ws.Range("E33:E39").PasteSpecial(-4122)

# Copy the activity-name formatting (style index 1) from an existing activity cell (C32) onto the new ones
$ws.Range("C32").Copy()
$ws.Range("C33:C39").PasteSpecial(-4122)

# Copy the trailing "I" column status-cell formatting (styles 8 / 24, alternating as in the source)
$ws.Range("I9").Copy()
$ws.Range("I33").PasteSpecial(-4122)
$ws.Range("I36").PasteSpecial(-4122)
$ws.Range("I38").PasteSpecial(-4122)

$ws.Range("J10").Copy()
$ws.Range("I34").PasteSpecial(-4122)
$ws.Range("I35").PasteSpecial(-4122)
$ws.Range("I37").PasteSpecial(-4122)
$ws.Range("I39").PasteSpecial(-4122)

# New activity index numbers for the added rows 37-39
$ws.Range("B37").Value = 31
$ws.Range("B38").Value = 32
$ws.Range("B39").Value = 33

# Activity names (rows 33-36 reuse already-existing text, 37-38 too, 39 is a brand-new string)
$ws.Range("C33").Value = "Leer archivo plano CSV"
$ws.Range("C34").Value = "Cargar CSV al vector de clientes"
$ws.Range("C35").Value = "Cambio en el total de clientes registrados JOptionPane"
$ws.Range("C36").Value = " Corrección de tabla para mostrar CSV"
$ws.Range("C37").Value = "Gráfico de columnas"
$ws.Range("C38").Value = "Crear archivo plano TXT"
$ws.Range("C39").Value = "Edición final de botones y orden burbuja para productos"

# Start / end dates for rows 33-39 (1 y 2 de septiembre 2022)
$ws.Range("D33:D39").Value = 44805
$ws.Range("E33:E39").Value = 44806

# --- Two closing "blank separator" rows (40-41) + one trailing row (42) ---
$ws.Range("C32").Copy()
$ws.Range("C40:C41").PasteSpecial(-4122)
$ws.Range("C42").PasteSpecial(-4122)

$ws.Range("J28").Copy()
$ws.Range("D40:I41").PasteSpecial(-4122)

# --- Remove the now-unused J28:J30 formatting-only cells ---
$ws.Range("J28:J30").Clear()

# --- Sheet view: scroll position, zoom and selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 2
$win.Zoom = 100
$ws.Range("G31").Select() | Out-Null
